# Apply the INDIANA_2022 data-cleaning fix:
#  1. Rename the header columns to snake_case machine-readable names.
#  2. Title-case the Spanish connector words (de/del/el/y/la/los/las) inside
#     state & municipality names (e.g. "Pabellón de Arteaga" -> "Pabellón De Arteaga").
#  3. Fix one mis-capitalised municipality name ("MonteMorelos" -> "Montemorelos").
#  4. Correct two floating point percentage values by 1 ULP.
#  5. Drop the trailing metadata/footnote rows (1245-1249) that don't belong
#     in the tabular data, and let the sheet dimension shrink accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Header renames ------------------------------------------------
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# --- 2. Capitalize Spanish connector words across every cell ----------
# Excel's Cells.Replace does a global, case-sensitive substring
# replacement by default, which is exactly what's needed here: every
# occurrence of " de " etc. (padded with spaces so we never touch a
# word that merely contains these letters, like "Madero" or "Cárdenas")
# becomes " De " and so on. Order does not matter since the patterns
# don't overlap (" de la " still leaves a bare " la " to match next).
$ws.Cells.Replace(" de ", " De ")
$ws.Cells.Replace(" del ", " Del ")
$ws.Cells.Replace(" el ", " El ")
$ws.Cells.Replace(" y ", " Y ")
$ws.Cells.Replace(" la ", " La ")
$ws.Cells.Replace(" los ", " Los ")
$ws.Cells.Replace(" las ", " Las ")

# --- 3. One-off name fix ------------------------------------------------
$ws.Cells.Replace("MonteMorelos", "Montemorelos")

# --- 4. Tiny floating point corrections --------------------------------
$ws.Range("D699").Value = 0.009766411598872333
$ws.Range("D847").Value = 0.009162303664921463

# --- 5. Remove the trailing metadata rows ------------------------------
$ws.Range("A1245:A1249").EntireRow.Delete()
